# Rename preview file and md file, move md jsonjs in db.
# Row 3 of the Tableau1 table (A3:D3) described a "folder_1_md" / "example_1"
# markdown example located in a subfolder; it is renamed/relocated to a
# top-level "tourisme_exemple" markdown example.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "tourisme_exemple"
$ws.Range("B3").Value = "md"
$ws.Range("C3").Value = "Tourisme exemple"
$ws.Range("D3").Value = "data/md/tourisme_exemple.md"

# Update the active selection to match the saved view state.
$ws.Range("D13").Select()
